$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10002, 110021),
    @(10003, 110022),
    @(10004, 110023),
    @(10005, 110024),
    @(10006, 110025),
    @(10007, 110026),
    @(10008, 110027),
    @(10009, 110028),
    @(10010, 110029)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

$ws.Range("A31:XFD1048576").Select() | Out-Null

$ws.PageSetup.Orientation = 1
